$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: "category_id" - copy formatting from the adjacent
# header cell E1 (same header style) then set its value.
$ws.Range("F1").Value = "category_id"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell F2: category_id value for the product row.
$ws.Range("F2").Value = 1

# Give the new column an explicit width (closest representable value to
# the bestFit width Excel would have computed for "category_id").
$ws.Columns("F").ColumnWidth = 9

# Update the active selection, matching the saved view state.
$ws.Range("F1").Select()
